$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 12: update title (D) and link (E)
$ws.Range("D12").Value = "TensorFlow 2.6.0 RC0"
$ws.Range("E12").Value = "https://tensorflow.blog/2021/06/30/tensorflow-2-6-0-rc0/"

# Row 51: update title (D) and link (E)
$ws.Range("D51").Value = "[python] 파이썬 예외 종류 정리(SyntaxError, TypeError, IndexError 등)"
$ws.Range("E51").Value = "https://bskyvision.com/1184"
